$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.92 = 42252.96 pesos`n✅ 42252.96 pesos = 9.91 = 972.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 100.845
$ws2.Range("O10").Value = 4261

$ws2.Range("N12").Value = 4264
$ws2.Range("O12").Value = 98.09999999999999
